$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.726.69"
$ws.Range("E2").Value = "  +5.83%  "

$ws.Range("D3").Value = "2.258.00"
$ws.Range("E3").Value = "  +4.45%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.22"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.35"
$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.414"
$ws.Range("E9").Value = "  +4.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "60.43"
$ws.Range("E10").Value = "  +4.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  +6.59%  "

$ws.Range("E12").Value = "  +1.73%  "

$ws.Range("D13").Value = "2.594.40"
$ws.Range("E13").Value = "  +4.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.24"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.70"
$ws.Range("E15").Value = "  +2.68%  "

$ws.Range("E16").Value = "  +2.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.69"
$ws.Range("E17").Value = "  +2.96%  "

$ws.Range("D18").Value = "2.275.90"
$ws.Range("E18").Value = "  +4.87%  "

$ws.Range("D19").Value = "41.579.24"
$ws.Range("E19").Value = "  +5.53%  "

$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  +9.94%  "

$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.75"
$ws.Range("E23").Value = "  +11.47%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +0.86%  "

$ws.Range("E27").Value = "  +5.90%  "

$ws.Range("E28").Value = "  +3.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.76"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("E30").Value = "  +3.34%  "

$ws.Range("E31").Value = "  +2.03%  "

$ws.Range("E32").Value = "  +8.24%  "

$ws.Range("E33").Value = "  +3.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.11"
$ws.Range("E34").Value = "  +8.10%  "

$ws.Range("E35").Value = "  +3.57%  "

$ws.Range("E36").Value = "  +4.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.94"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.86"
$ws.Range("E38").Value = "  +8.10%  "

$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000260"
$ws.Range("E40").Value = "  +64.77%  "

$ws.Range("E41").Value = "  +20.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  +6.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.81"
$ws.Range("E44").Value = "  +12.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.98"
$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0999"
$ws.Range("E46").Value = "  +7.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.68"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.24"
$ws.Range("E48").Value = "  +3.14%  "

$ws.Range("D49").Value = "1.512.60"
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("E50").Value = "  +3.54%  "

$ws.Range("E51").Value = "  -0.95%  "
